# Update column A values on the active sheet ("main").
# The diff replaces every existing value in A1:A90 with a new value, and
# appends new values for rows 91-150 (column A now spans A1:A150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    73.65148208090005, 71.5654190831555, 12.15103083423297, 11.426392493482673, 6.3657058416708665,
    6.2668354387146215, 6.018524849238329, 5.649278881567422, 4.894450650526875, 4.846680506968446,
    3.3180148797316704, 3.1468888293667807, 2.9178419344113253, 2.8721289066602838, 2.670463918754725,
    2.630459087462568, 2.2475412465946856, 2.144898641005488, 1.9074099873174817, 1.8695986000392686,
    1.6913328205314564, 1.5945974269300172, 1.3933967291190794, 1.3856417119495559, 1.356421292285049,
    1.3353592588978371, 1.3134652563017128, 1.3050398746224352, 1.2230034804802028, 1.2184526866887373,
    1.1370577909168755, 1.1166461095520606, 0.9126170840170339, 0.7573962970162427, 0.5061659406008099,
    0.44811322979533025, 0.299233302410026, 0.23882809096418925, 0.14763480396613896, 0.10690207544706783,
    0.08066015286641881, 0.06117751155458824, 0.0516237354924994, 0.039216235553200175, 0.033924354672567075,
    0.02818353484368433, 0.018706762401778276, 0.01606118750660833, 0.01536831940283495, 0.009854138697556534,
    0.0090437283024589, 0.008995381662573887, 0.006607137195720241, 0.0056186561312976425, 0.005459864843209241,
    0.00467372986271009, 0.004147896918771631, 0.0038359442959082692, 0.003657711172858614, 0.003455151786063285,
    0.003385924633234284, 0.0032070780191149647, 0.003076779081074718, 0.0030534422443745856, 0.0027970848477307116,
    0.002717332598026908, 0.002699221848845295, 0.0025683882405719847, 0.0024886911039145596, 0.0023960851416382324,
    0.0023680908144621668, 0.0023045919865055924, 0.0020452421299154743, 0.0020069690815275555, 0.002006822961407528,
    0.0019843201841516995, 0.001966889910796623, 0.0018018616106045399, 0.0017809898157123698, 0.0017736190299969129,
    0.001662841246094959, 0.0015985812173693512, 0.0015582334425136564, 0.0015268716892666548, 0.001469457763706943,
    0.001467536430491016, 0.0014037750422028802, 0.0012506835340837775, 0.0012367690572032172, 0.0012241624914952744,
    0.001204129189837505, 0.0011471813705538677, 0.0011395248216345575, 0.0011357566921070962, 0.0010977174281391901,
    0.0010866867461262063, 0.0010806234782664691, 0.0010405661219855275, 0.0009966494655046908, 0.0009949099246063606,
    0.0009787861796821961, 0.000977803518505653, 0.0009743933140480024, 0.0009732708801080717, 0.000940511097573523,
    0.0009118514251224665, 0.0008744708249892553, 0.0008526734599615674, 0.0008219158188079482, 0.0007878001323665303,
    0.0007728887291290534, 0.0007689752150222075, 0.0007518125389021705, 0.0007377389184933806, 0.0007167198878757198,
    0.0007050346412913177, 0.000703134373748179, 0.0006810349004230626, 0.0006778373291260104, 0.0006548477777994086,
    0.0006508711227124996, 0.0006486590664545881, 0.0006483967356305447, 0.0006406891664693552, 0.0006363472058167124,
    0.0006082172221525704, 0.0006069996074795361, 0.0005804914076622242, 0.0005781974380133214, 0.0005391348426972122,
    0.0004824044080269323, 0.0004797384102598265, 0.00046712699296051036, 0.00045795791303395193, 0.00045181829095976083,
    0.00045103939804816477, 0.0004445405538737733, 0.00041448600340059186, 0.0004128167317757558, 0.00040628564781553206,
    0.00040309505131717335, 0.0003980364371124916, 0.0003509577396593601, 0.0003280882994556986, 0.0003076011747652274,
    0.0003039851943664237, 0.000290357114979041, 0.0002862857134435256, 0.00022244058040534642, 0.00021445972501369243
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
